# Update Name of Algo
# Applies corrected KNN imputation result values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.269999999999999
$ws.Range("A8").Value = -21.107
$ws.Range("A10").Value = -20.945
$ws.Range("A12").Value = -21.694
$ws.Range("C13").Value = -13.059
$ws.Range("A18").Value = -21.694
$ws.Range("D20").Value = -8.222
